# Journal_Travail_Daniel.xlsx - add new journal entries for work done on
# "Updated IOTransaction for shared budgets"
#
# The "Feuil1" worksheet lists journal entries in rows 5-75 (columns A=Date,
# B=Activité, C=Heures), with a running total in C76. Rows 66-68 were blank
# placeholder rows; we fill them with three new entries dated 43236
# (2018-05-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 66
$ws.Cells.Item(66, 1).Value = 43236
$ws.Cells.Item(66, 2).Value = "Implémentation de quelques fonctions demandées par la GUI pour les dettes / budgets"
$ws.Cells.Item(66, 3).Value = 1.25

# Row 67
$ws.Cells.Item(67, 1).Value = 43236
$ws.Cells.Item(67, 2).Value = "Merge de la branche master dans la branch fb-derby pour mettre à jour et implémenter la fin de derby"
$ws.Cells.Item(67, 3).Value = 1.5

# Row 68
$ws.Cells.Item(68, 1).Value = 43236
$ws.Cells.Item(68, 2).Value = "Corrections dans quelques problèmes dans la bll"
$ws.Cells.Item(68, 3).Value = 1.5

# Rows 66 and 67 wrap onto two lines in the original sheet (ht="30").
$ws.Rows.Item(66).RowHeight = 30
$ws.Rows.Item(67).RowHeight = 30

# Move the active selection the way the author's session ended up (on A69).
$ws.Range("A69").Select()
